$d = $word.ActiveDocument
$d.Content.Find.Execute("Implementations", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Implemeations", 2)
